$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B5 text (row for 2012-07-23, was "Anlegen SVN und Verzeichnisstruktur")
$ws.Range("B5").Value = "Anlegen SVN und Verzeichnisstruktur, Create project sort, Implement test for own quicksort, c qsort and c++ sort"

# Add new row 6: date 2012-07-29, activity text
$ws.Range("A6").NumberFormat = "mm-dd-yy"
$ws.Range("A6").Value = (Get-Date -Year 2012 -Month 7 -Day 29).Date
$ws.Range("B6").Value = "Created OpenCL classes, Implemented Bealto ParallelSelectionSort"

# Update selection to B15
$ws.Range("B15").Select()
